$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 8).Value = "none"
}

$ws.Range("J20").Select()
